$d = $word.ActiveDocument

# 1. Fix typo: "ürlap" -> "űrlap" (both occurrences)
$d.Content.Find.Execute("ürlap", $false, $false, $false, $false, $false, $true, 1, $false, "űrlap", 2)

Write-Output "done"
